$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - set values then copy the header formatting (bold, centered,
# bordered) from the existing H1 header cell so I1/J1 match the other headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-23 for columns I and J
$data = @(
    @(2, 6, 7),
    @(3, 8, 8),
    @(4, 7, 7),
    @(5, 9, 9),
    @(6, 7, 7),
    @(7, 6, 7),
    @(8, 5, 7),
    @(9, 8, 9),
    @(10, 8, 9),
    @(11, 8, 9),
    @(12, 8, 9),
    @(13, 1, 4),
    @(14, 1, 4),
    @(15, 1, 5),
    @(16, 1, 6),
    @(17, 1, 5),
    @(18, 1, 6),
    @(19, 1, 5),
    @(20, 7, 9),
    @(21, 1, 3),
    @(22, 7, 7),
    @(23, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
